$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit inserts one new record as row 35, pushing every subsequent
# existing record (old rows 35..145) down by one row (new rows 36..146).
# Work from the bottom up so we never overwrite data before reading it.
for ($r = 146; $r -ge 36; $r--) {
    for ($c = 1; $c -le 18; $c++) {
        $v = $ws.Cells.Item($r - 1, $c).Value2
        $ws.Cells.Item($r, $c).Value2 = $v
    }
}

# Make sure the newly written last row (146) keeps the same date number
# format as column D uses everywhere else in the table.
$ws.Cells.Item(146, 4).NumberFormat = $ws.Cells.Item(145, 4).NumberFormat

# Now overwrite row 35 with the new record's data (same as old row 35
# except for the Fecha, Volumen, Precio promedio ponderado and Precio
# $/Kg columns).
$ws.Cells.Item(35, 1).Value2 = 11
$ws.Cells.Item(35, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(35, 3).Value2 = "Bíobío"
$ws.Cells.Item(35, 4).Value2 = 44623
$ws.Cells.Item(35, 5).Value2 = 8
$ws.Cells.Item(35, 6).Value2 = 100112003
$ws.Cells.Item(35, 7).Value2 = "Ajo"
$ws.Cells.Item(35, 8).Value2 = "Chino"
$ws.Cells.Item(35, 9).Value2 = "Primera"
$ws.Cells.Item(35, 10).Value2 = 220
$ws.Cells.Item(35, 11).Value2 = 16000
$ws.Cells.Item(35, 12).Value2 = 17000
$ws.Cells.Item(35, 13).Value2 = 16455
$ws.Cells.Item(35, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(35, 15).Value2 = "China"
$ws.Cells.Item(35, 16).Value2 = 1646
$ws.Cells.Item(35, 17).Value2 = 10
$ws.Cells.Item(35, 18).Value2 = "Hortaliza"
